$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The TODO item "If a non-SBT base type is deleted, either the type(s) that
# point to it have to be nulled, or better, it should be prohibited in the
# UI." (previously tracked in the TO DO list) has now been implemented.

# 1. Mark it DONE: fill the previously-blank spacer row (A36) in the DONE
#    section with that note text, matching the DONE section's formatting
#    (strikethrough + wrap, copied from a neighboring DONE row) and let the
#    row grow tall enough to show the full, wrapped note.
$ws.Rows("36").ClearFormats()
$ws.Range("A36").Value = "If a non-SBT base type is deleted, either the type(s) that point to it have to be nulled, or better, it should be prohibited in the UI."
$ws.Range("A37").Copy()
$ws.Range("A36").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B36").Clear()
$ws.Rows("36").RowHeight = 43.5

# 2. Remove the corresponding row from the TO DO list (old row 47), which
#    shifts the remaining rows (old row 48, "Add usergroups") up by one.
$ws.Rows("47").Delete()

$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A46").Select()
